$wb = $excel.ActiveWorkbook

# --- Update "BAU Emissions" sheet ---
$wsBau = $wb.Worksheets.Item("BAU Emissions")

# 1) Replace the " : NoSettings" suffix with " : test" across all the
#    descriptive labels in column A (these are the shared strings that
#    changed across ~275 rows).
$wsBau.Range("A1:A300").Replace(" : NoSettings", " : test") | Out-Null

# 2) Update the data values for row 94 (columns M through AE).
$cols = @("M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE")
$vals = @(1001080,2002150,3003230,4004300,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380,5005380)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsBau.Range("$($cols[$i])94").Value2 = $vals[$i]
}

# 3) Update the view/selection on this sheet.
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select() | Out-Null

# --- Update "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value2 = 45387

# --- Make "About" the active/selected sheet (tabSelected moves here,
#     and is cleared from "Current and Planned Capacity"). ---
$wsAbout.Activate()
